# Insert 5 new data rows for "Vega Modelo de Temuco - Repollo" (weekly update),
# pushing the existing rows 1212-1277 down to 1217-1282, then populate the
# 5 newly inserted rows (1212-1216) with the new weekly price records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows before the current row 1212 (shifts old 1212:1277 -> 1217:1282)
$ws.Rows.Item(1212).Resize(5).Insert()

# Constant columns shared by every data row in this sheet
$mercadoId = 10
$mercado   = "Vega Modelo de Temuco"
$region    = "La Araucanía"
$codreg    = 9
$catId     = 100112006
$categoria = "Repollo"
$calidad   = "Primera"
$unidad    = "$/unidad"
$kgUnidad  = 1
$clasif    = "Hortaliza"

function Set-DataRow {
    param($row, $fecha, $variedad, $volumen, $precioMin, $precioMax, $precioProm, $origen, $precioKg)

    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $catId
    $ws.Cells.Item($row, 7).Value = $categoria
    $ws.Cells.Item($row, 8).Value = $variedad
    $ws.Cells.Item($row, 9).Value = $calidad
    $ws.Cells.Item($row, 10).Value = $volumen
    $ws.Cells.Item($row, 11).Value = $precioMin
    $ws.Cells.Item($row, 12).Value = $precioMax
    $ws.Cells.Item($row, 13).Value = $precioProm
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $precioKg
    $ws.Cells.Item($row, 17).Value = $kgUnidad
    $ws.Cells.Item($row, 18).Value = $clasif
}

Set-DataRow 1212 45267 "Crespo record" 1500 1400 1400 1400 "Provincia del Elquí" 1400
Set-DataRow 1213 45267 "Crespo record" 3000 1400 1400 1400 "Región Metropolitana" 1400
Set-DataRow 1214 45267 "Crespo record" 5000 1300 1500 1420 "Región del Maule" 1420
Set-DataRow 1215 45267 "Morada(o)"     600  1500 1500 1500 "Provincia del Elquí" 1500
Set-DataRow 1216 45267 "Morada(o)"     800  1400 1400 1400 "Región del Maule" 1400
